$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.85519579756308
$ws.Range("C2").Value = 14.97343076181612
$ws.Range("E2").Value = 16.49170269138759
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.618659601243766
$ws.Range("I2").Value = 27.50212814490051

$ws.Range("B3").Value = 15.18423979154654
$ws.Range("C3").Value = 14.08573814935257
$ws.Range("E3").Value = 15.54708399255748
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.624069924259701
$ws.Range("I3").Value = 26.92483206538295

$ws.Range("B4").Value = 14.7635238801976
$ws.Range("C4").Value = 13.51619194847083
$ws.Range("E4").Value = 14.94359293374112
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.627547748268274
$ws.Range("I4").Value = 26.5709284640591

$ws.Range("B5").Value = 14.59018781463931
$ws.Range("C5").Value = 13.27815271822352
$ws.Range("E5").Value = 14.6920276531438
$ws.Range("F5").Value = 15.00819731993643
$ws.Range("G5").Value = 3.629004409140012
$ws.Range("I5").Value = 26.42700335738451

$ws.Range("B6").Value = 14.56130097741652
$ws.Range("C6").Value = 13.23827468808588
$ws.Range("E6").Value = 14.6499238518041
$ws.Range("F6").Value = 14.96433081551586
$ws.Range("G6").Value = 3.629248674001613
$ws.Range("I6").Value = 26.4031269174809

$ws.Range("B7").Value = 14.76119343266333
$ws.Range("C7").Value = 13.51300541770711
$ws.Range("E7").Value = 14.94022266838935
$ws.Range("F7").Value = 15.26647399323726
$ws.Range("G7").Value = 3.627567233398814
$ws.Range("I7").Value = 26.56898604855467

$ws.Range("B8").Value = 15.62583981500083
$ws.Range("C8").Value = 14.67255502284713
$ws.Range("E8").Value = 16.17101947384857
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.620492867985571
$ws.Range("I8").Value = 27.30305192504547

$ws.Range("B9").Value = 17.24019302504579
$ws.Range("C9").Value = 16.74466533879647
$ws.Range("E9").Value = 18.48074792959202
$ws.Range("F9").Value = 19.0027458068253
$ws.Range("G9").Value = 3.607846383200301
$ws.Range("I9").Value = 28.74036819957285

$ws.Range("B10").Value = 18.36321100165532
$ws.Range("C10").Value = 18.13658108011846
$ws.Range("E10").Value = 20.15464664419617
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.599287963418478
$ws.Range("I10").Value = 29.78643701131075

$ws.Range("B11").Value = 18.85815080831965
$ws.Range("C11").Value = 18.74044185307907
$ws.Range("E11").Value = 20.87509341356871
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.595550483900027
$ws.Range("I11").Value = 30.25844663507181

$ws.Range("B12").Value = 19.04312845025527
$ws.Range("C12").Value = 18.96483140483337
$ws.Range("E12").Value = 21.14205588196877
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.594157347008706
$ws.Range("I12").Value = 30.43649549188559

$ws.Range("B13").Value = 19.00340118397767
$ws.Range("C13").Value = 18.9166962424568
$ws.Range("E13").Value = 21.08482045518352
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.594456401940157
$ws.Range("I13").Value = 30.39818231475626

$ws.Range("B14").Value = 18.8734188435355
$ws.Range("C14").Value = 18.75898857927392
$ws.Range("E14").Value = 20.89717367301296
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.595435426824915
$ws.Range("I14").Value = 30.27310942679796

$ws.Range("B15").Value = 18.79347825092647
$ws.Range("C15").Value = 18.66182921719972
$ws.Range("E15").Value = 20.78147338103637
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.596037987369376
$ws.Range("I15").Value = 30.19640475423341

$ws.Range("B16").Value = 18.33053123681633
$ws.Range("C16").Value = 18.09652150633615
$ws.Range("E16").Value = 20.10674206141509
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.599535331175743
$ws.Range("I16").Value = 29.7555006725853

$ws.Range("B17").Value = 18.04233106977391
$ws.Range("C17").Value = 17.7421621922296
$ws.Range("E17").Value = 19.68234129455281
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.601720571798645
$ws.Range("I17").Value = 29.48393380272258

$ws.Range("B18").Value = 17.87507452982459
$ws.Range("C18").Value = 17.53558604760779
$ws.Range("E18").Value = 19.43437378506921
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.602992142262557
$ws.Range("I18").Value = 29.32738086253836

$ws.Range("B19").Value = 17.81819330150951
$ws.Range("C19").Value = 17.46517138649667
$ws.Range("E19").Value = 19.34975138532903
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.60342520227751
$ws.Range("I19").Value = 29.27431807573767

$ws.Range("B20").Value = 18.07316602484755
$ws.Range("C20").Value = 17.78017029473189
$ws.Range("E20").Value = 19.72791901584559
$ws.Range("F20").Value = 20.2495528364879
$ws.Range("G20").Value = 3.601486431810198
$ws.Range("I20").Value = 29.51288030243327

$ws.Range("B21").Value = 18.91166529336682
$ws.Range("C21").Value = 18.80542769788019
$ws.Range("E21").Value = 20.95244864189617
$ws.Range("F21").Value = 21.46857628470567
$ws.Range("G21").Value = 3.595147263817091
$ws.Range("I21").Value = 30.30986617177847

$ws.Range("B22").Value = 19.44536384345048
$ws.Range("C22").Value = 19.46282356137354
$ws.Range("E22").Value = 21.71864173808762
$ws.Range("F22").Value = 22.22866616901555
$ws.Range("G22").Value = 3.591133336931493
$ws.Range("I22").Value = 30.82665710157598

$ws.Range("B23").Value = 19.16187349425923
$ws.Range("C23").Value = 19.10852766094444
$ws.Range("E23").Value = 21.31281661684445
$ws.Range("F23").Value = 21.82633154475864
$ws.Range("G23").Value = 3.59326391271495
$ws.Range("I23").Value = 30.55125313582716

$ws.Range("B24").Value = 18.05923040870989
$ws.Range("C24").Value = 17.76299570462717
$ws.Range("E24").Value = 19.70732570899519
$ws.Range("F24").Value = 20.22900810905294
$ws.Range("G24").Value = 3.601592238990905
$ws.Range("I24").Value = 29.49979490486466

$ws.Range("B25").Value = 16.81367981890923
$ws.Range("C25").Value = 16.2065894513379
$ws.Range("E25").Value = 17.82765916562611
$ws.Range("F25").Value = 18.34778573295697
$ws.Range("G25").Value = 3.611137813282459
$ws.Range("I25").Value = 28.35267522503778
